{"js": "// Replace the three-digit / one-digit division problems throughout the\n// document body (including inside the table cells) with the new values\n// from the commit. Each old value is unique within the document, so a\n// simple search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"116\u00f78=\", \"705\u00f74=\"],\n  [\"852\u00f78=\", \"372\u00f79=\"],\n  [\"526\u00f73=\", \"796\u00f77=\"],\n  [\"971\u00f76=\", \"289\u00f73=\"],\n  [\"134\u00f79=\", \"346\u00f79=\"],\n  [\"346\u00f73=\", \"297\u00f78=\"],\n  [\"464\u00f78=\", \"304\u00f72=\"],\n  [\"509\u00f77=\", \"607\u00f79=\"],\n  [\"404\u00f78=\", \"734\u00f77=\"],\n  [\"915\u00f79=\", \"225\u00f77=\"],\n  [\"605\u00f72=\", \"345\u00f74=\"],\n  [\"601\u00f78=\", \"921\u00f74=\"],\n  [\"165\u00f72=\", \"131\u00f79=\"],\n  [\"691\u00f79=\", \"269\u00f72=\"],\n  [\"475\u00f79=\", \"486\u00f77=\"],\n  [\"629\u00f76=\", \"969\u00f79=\"],\n  [\"180\u00f76=\", \"566\u00f74=\"],\n  [\"596\u00f75=\", \"513\u00f73=\"],\n  [\"986\u00f78=\", \"819\u00f79=\"],\n  [\"911\u00f75=\", \"565\u00f77=\"],\n  [\"952\u00f72=\", \"835\u00f78=\"],\n  [\"653\u00f75=\", \"440\u00f72=\"],\n  [\"104\u00f73=\", \"329\u00f74=\"],\n  [\"294\u00f78=\", \"914\u00f74=\"],\n  [\"155\u00f72=\", \"606\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit / one-digit division problems throughout the\n# document (including inside the table cells) with the new values from\n# the commit. Each old value is unique within the document, so a plain\n# Find/Replace (wdReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"116\u00f78=\", \"705\u00f74=\"),\n    @(\"852\u00f78=\", \"372\u00f79=\"),\n    @(\"526\u00f73=\", \"796\u00f77=\"),\n    @(\"971\u00f76=\", \"289\u00f73=\"),\n    @(\"134\u00f79=\", \"346\u00f79=\"),\n    @(\"346\u00f73=\", \"297\u00f78=\"),\n    @(\"464\u00f78=\", \"304\u00f72=\"),\n    @(\"509\u00f77=\", \"607\u00f79=\"),\n    @(\"404\u00f78=\", \"734\u00f77=\"),\n    @(\"915\u00f79=\", \"225\u00f77=\"),\n    @(\"605\u00f72=\", \"345\u00f74=\"),\n    @(\"601\u00f78=\", \"921\u00f74=\"),\n    @(\"165\u00f72=\", \"131\u00f79=\"),\n    @(\"691\u00f79=\", \"269\u00f72=\"),\n    @(\"475\u00f79=\", \"486\u00f77=\"),\n    @(\"629\u00f76=\", \"969\u00f79=\"),\n    @(\"180\u00f76=\", \"566\u00f74=\"),\n    @(\"596\u00f75=\", \"513\u00f73=\"),\n    @(\"986\u00f78=\", \"819\u00f79=\"),\n    @(\"911\u00f75=\", \"565\u00f77=\"),\n    @(\"952\u00f72=\", \"835\u00f78=\"),\n    @(\"653\u00f75=\", \"440\u00f72=\"),\n    @(\"104\u00f73=\", \"329\u00f74=\"),\n    @(\"294\u00f78=\", \"914\u00f74=\"),\n    @(\"155\u00f72=\", \"606\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
